$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.491755872964859
$ws.Range("B1").Value = 2.56944751739502
$ws.Range("C1").Value = 6.328853130340576
$ws.Range("D1").Value = 1.54465115070343
$ws.Range("E1").Value = 0.8915167450904846
